$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-25 08:02:57"
$wsZhCn.Range("H2").Value = "2016-03-25 08:03:44"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-25 08:03:02"
$wsDeDe.Range("H2").Value = "2016-03-25 08:03:52"
